$wb = $excel.ActiveWorkbook

# diff hunk @ 923 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 295.7143
$ws.Range("I6").Value = 8
$ws.Range("K6").Value = 24
$ws.Range("M6").Value = 88

# diff hunk @ 2292 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1907.3334
$ws.Range("I34").Value = 1907.3334
$ws.Range("K34").Value = 1907.3334
$ws.Range("M34").Value = -1704.3334

# diff hunk @ 2393 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 1907.3334
$ws.Range("I36").Value = 1907.3334
$ws.Range("K36").Value = 1907.3334
$ws.Range("M36").Value = -1192.3334

# diff hunk @ 3762 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2735.25
$ws.Range("I64").Value = 2797.3333
$ws.Range("K64").Value = 2797.3333
$ws.Range("M64").Value = -2549.3333

# diff hunk @ 3912 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2735.25
$ws.Range("I67").Value = 2797.3333
$ws.Range("K67").Value = 2797.3333
$ws.Range("M67").Value = -1939.3333

# diff hunk @ 4270 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9932
$ws.Range("I74").Value = 9932
$ws.Range("K74").Value = 9932
$ws.Range("M74").Value = -8996

# diff hunk @ 4365 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5333.3335
$ws.Range("I76").Value = 4000
$ws.Range("K76").Value = 4000
$ws.Range("M76").Value = -3685

# diff hunk @ 4414 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 9932
$ws.Range("I77").Value = 9932
$ws.Range("K77").Value = 49660
$ws.Range("M77").Value = -44980

# diff hunk @ 4509 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5333.3335
$ws.Range("I79").Value = 4000
$ws.Range("K79").Value = 4000
$ws.Range("M79").Value = -2908

# diff hunk @ 4558 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1632.25
$ws.Range("I80").Value = 1379
$ws.Range("K80").Value = 4137
$ws.Range("M80").Value = -3139

# diff hunk @ 4705 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1632.25
$ws.Range("I83").Value = 1379
$ws.Range("K83").Value = 12411
$ws.Range("M83").Value = -7419

# diff hunk @ 6766 -> sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 629.5
$ws.Range("I125").Value = 252
$ws.Range("J125").Value = 663.8182
$ws.Range("K125").Value = 2268
$ws.Range("L125").Value = 5974.3638
$ws.Range("M125").Value = 192
$ws.Range("N125").Value = -10894.3638

# diff hunk @ 9782 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2524.3572
$ws.Range("I45").Value = 1631.909
$ws.Range("K45").Value = 1631.909
$ws.Range("M45").Value = -1254.909

# diff hunk @ 10548 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1829.6
$ws.Range("I61").Value = 1712
$ws.Range("K61").Value = 1712
$ws.Range("M61").Value = -1500

# diff hunk @ 10646 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5092.4287
$ws.Range("I63").Value = 1229.4
$ws.Range("K63").Value = 1229.4
$ws.Range("M63").Value = -543.4000000000001

# diff hunk @ 10790 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5092.4287
$ws.Range("I66").Value = 1229.4
$ws.Range("K66").Value = 6147
$ws.Range("M66").Value = -2715

# diff hunk @ 13949 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6790.4287
$ws.Range("I132").Value = 6061.7036
$ws.Range("K132").Value = 18185.1108
$ws.Range("M132").Value = -15655.1108

# diff hunk @ 14139 -> sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1829.6
$ws.Range("I136").Value = 1712
$ws.Range("K136").Value = 5136
$ws.Range("M136").Value = -2586

# diff hunk @ 15403 -> sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2295.4614
$ws.Range("I20").Value = 2226.4546
$ws.Range("K20").Value = 2226.4546
$ws.Range("M20").Value = -1979.4546

# diff hunk @ 19585 -> sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1870.7
$ws.Range("I107").Value = 1634.2222
$ws.Range("K107").Value = 1634.2222
$ws.Range("M107").Value = 285.7778000000001

# diff hunk @ 22853 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# diff hunk @ 25847 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 18915
$ws.Range("J96").Value = 18915
$ws.Range("L96").Value = 18915
$ws.Range("N96").Value = -24407

# diff hunk @ 26380 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 834.35
$ws.Range("I107").Value = 699.2941
$ws.Range("K107").Value = 699.2941
$ws.Range("M107").Value = 1220.7059

# diff hunk @ 27676 -> sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2213.1667
$ws.Range("I134").Value = 1916
$ws.Range("K134").Value = 5748
$ws.Range("M134").Value = -3213

# diff hunk @ 28973 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 439.14285
$ws.Range("I18").Value = 354
$ws.Range("K18").Value = 1062
$ws.Range("M18").Value = -893

# diff hunk @ 29221 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 978
$ws.Range("I23").Value = 900
$ws.Range("J23").Value = 997.5
$ws.Range("K23").Value = 2700
$ws.Range("L23").Value = 2992.5
$ws.Range("M23").Value = -2465
$ws.Range("N23").Value = -3462.5

# diff hunk @ 30345 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1837.5
$ws.Range("I46").Value = 1300
$ws.Range("J46").Value = 2375
$ws.Range("K46").Value = 3900
$ws.Range("L46").Value = 7125
$ws.Range("M46").Value = -3809
$ws.Range("N46").Value = -7307

# diff hunk @ 30847 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9147.429
$ws.Range("I56").Value = 9147.429
$ws.Range("K56").Value = 9147.429
$ws.Range("M56").Value = -8617.429

# diff hunk @ 32081 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1475
$ws.Range("I81").Value = 1475
$ws.Range("K81").Value = 4425
$ws.Range("M81").Value = -3302

# diff hunk @ 32225 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 1475
$ws.Range("I84").Value = 1475
$ws.Range("K84").Value = 13275
$ws.Range("M84").Value = -7659

# diff hunk @ 33643 -> sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1122.0189
$ws.Range("I113").Value = 1136.4166
$ws.Range("J113").Value = 983.8
$ws.Range("K113").Value = 3409.2498
$ws.Range("L113").Value = 2951.4
$ws.Range("M113").Value = -1239.2498
$ws.Range("N113").Value = -7291.4

# diff hunk @ 37829 -> sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 35998
$ws.Range("J57").Value = 35998
$ws.Range("L57").Value = 35998
$ws.Range("N57").Value = -37638

# diff hunk @ 39956 -> sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 880.6
$ws.Range("I102").Value = 880.6
$ws.Range("K102").Value = 880.6
$ws.Range("M102").Value = 741.4

# diff hunk @ 40480 -> sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2764.7334
$ws.Range("I113").Value = 2773.3333
$ws.Range("J113").Value = 2751.8333
$ws.Range("K113").Value = 2773.3333
$ws.Range("L113").Value = 2751.8333
$ws.Range("M113").Value = -603.3332999999998
$ws.Range("N113").Value = -7091.8333

# diff hunk @ 40906 -> sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2481.25
$ws.Range("I122").Value = 2481.25
$ws.Range("K122").Value = 7443.75
$ws.Range("M122").Value = -4993.75

# diff hunk @ 42636 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 545.2308
$ws.Range("I16").Value = 659.7
$ws.Range("J16").Value = 163.66667
$ws.Range("K16").Value = 659.7
$ws.Range("L16").Value = 163.66667
$ws.Range("M16").Value = -489.7
$ws.Range("N16").Value = -503.66667

# diff hunk @ 45154 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9459.799999999999
$ws.Range("J68").Value = 9950
$ws.Range("L68").Value = 9950
$ws.Range("N68").Value = -11448

# diff hunk @ 45301 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 9459.799999999999
$ws.Range("J71").Value = 9950
$ws.Range("L71").Value = 49750
$ws.Range("N71").Value = -57238

# diff hunk @ 46683 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1539.8
$ws.Range("I100").Value = 1539.8
$ws.Range("K100").Value = 1539.8
$ws.Range("M100").Value = -998.8

# diff hunk @ 48209 -> sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3996.3333
$ws.Range("I132").Value = 3996.3333
$ws.Range("K132").Value = 11988.9999
$ws.Range("M132").Value = -9458.999899999999

# diff hunk @ 51621 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 100000
$ws.Range("I61").Value = 100000
$ws.Range("K61").Value = 100000
$ws.Range("M61").Value = -99708

# diff hunk @ 51670 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8799.75
$ws.Range("I62").Value = 8799.75
$ws.Range("K62").Value = 8799.75
$ws.Range("M62").Value = -8175.75

# diff hunk @ 51814 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8799.75
$ws.Range("I65").Value = 8799.75
$ws.Range("K65").Value = 43998.75
$ws.Range("M65").Value = -40878.75

# diff hunk @ 54124 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 466.69232
$ws.Range("I113").Value = 400.55554
$ws.Range("K113").Value = 1201.66662
$ws.Range("M113").Value = 968.33338

# diff hunk @ 54761 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530

# diff hunk @ 55046 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4224.25
$ws.Range("I132").Value = 4224.25
$ws.Range("K132").Value = 12672.75
$ws.Range("M132").Value = -10142.75

# diff hunk @ 55481 -> sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 78125.25
$ws.Range("J141").Value = 70833.336
$ws.Range("L141").Value = 70833.336
$ws.Range("N141").Value = -81193.336
